$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header swap: a new data column was added, shifting the existing
# "average_doctor" header (and its values) into the "_old" slot, and
# recomputing fresh averages for the new harvard-case-classification pass.
$ws.Range("BP1").Value = "average_doctor_old"
$ws.Range("BQ1").Value = "average_doctor"

# --- Updated statistics values (rows 4-13) ---

# row 4
$ws.Cells.Item(4, 5).Value = 0.423
$ws.Cells.Item(4, 6).Value = 0.07199999999999999
$ws.Cells.Item(4, 7).Value = 0.268
$ws.Cells.Item(4, 14).Value = 0.432
$ws.Cells.Item(4, 15).Value = 0.064
$ws.Cells.Item(4, 16).Value = 0.254
$ws.Cells.Item(4, 17).Value = 0.023
$ws.Cells.Item(4, 18).Value = 0.016
$ws.Cells.Item(4, 19).Value = 0.128
$ws.Cells.Item(4, 23).Value = 0.283
$ws.Cells.Item(4, 24).Value = 0.109
$ws.Cells.Item(4, 25).Value = 0.33
$ws.Cells.Item(4, 35).Value = 0.288
$ws.Cells.Item(4, 36).Value = 0.08799999999999999
$ws.Cells.Item(4, 37).Value = 0.297
$ws.Cells.Item(4, 47).Value = 0.188
$ws.Cells.Item(4, 48).Value = 0.028
$ws.Cells.Item(4, 49).Value = 0.168
$ws.Cells.Item(4, 53).Value = 2.012
$ws.Cells.Item(4, 54).Value = 0.158
$ws.Cells.Item(4, 55).Value = 0.397
$ws.Cells.Item(4, 59).Value = 0.732
$ws.Cells.Item(4, 60).Value = 0.141
$ws.Cells.Item(4, 61).Value = 0.375
$ws.Cells.Item(4, 65).Value = 0.719
$ws.Cells.Item(4, 66).Value = 0.08
$ws.Cells.Item(4, 67).Value = 0.282
$ws.Cells.Item(4, 68).Value = 0.671
$ws.Cells.Item(4, 69).Value = 0.71

# row 5
$ws.Cells.Item(5, 5).Value = 0.538
$ws.Cells.Item(5, 6).Value = 0.08500000000000001
$ws.Cells.Item(5, 7).Value = 0.292
$ws.Cells.Item(5, 14).Value = 0.733
$ws.Cells.Item(5, 15).Value = 0.081
$ws.Cells.Item(5, 16).Value = 0.285
$ws.Cells.Item(5, 17).Value = 0.015
$ws.Cells.Item(5, 18).Value = 0.007
$ws.Cells.Item(5, 19).Value = 0.082
$ws.Cells.Item(5, 23).Value = 0.273
$ws.Cells.Item(5, 24).Value = 0.109
$ws.Cells.Item(5, 25).Value = 0.33
$ws.Cells.Item(5, 35).Value = 0.309
$ws.Cells.Item(5, 36).Value = 0.098
$ws.Cells.Item(5, 37).Value = 0.314
$ws.Cells.Item(5, 47).Value = 0.365
$ws.Cells.Item(5, 48).Value = 0.093
$ws.Cells.Item(5, 49).Value = 0.306
$ws.Cells.Item(5, 53).Value = 1.337
$ws.Cells.Item(5, 54).Value = 0.079
$ws.Cells.Item(5, 55).Value = 0.28
$ws.Cells.Item(5, 59).Value = 0.398
$ws.Cells.Item(5, 60).Value = 0.051
$ws.Cells.Item(5, 61).Value = 0.225
$ws.Cells.Item(5, 65).Value = 0.549
$ws.Cells.Item(5, 66).Value = 0.063
$ws.Cells.Item(5, 67).Value = 0.252
$ws.Cells.Item(5, 68).Value = 0.446
$ws.Cells.Item(5, 69).Value = 0.456

# row 6
$ws.Cells.Item(6, 5).Value = 0.474
$ws.Cells.Item(6, 14).Value = 0.544
$ws.Cells.Item(6, 17).Value = 0.018
$ws.Cells.Item(6, 23).Value = 0.278
$ws.Cells.Item(6, 35).Value = 0.298
$ws.Cells.Item(6, 47).Value = 0.248
$ws.Cells.Item(6, 53).Value = 1.599
$ws.Cells.Item(6, 59).Value = 0.516
$ws.Cells.Item(6, 65).Value = 0.623
$ws.Cells.Item(6, 68).Value = 0.533
$ws.Cells.Item(6, 69).Value = 0.552

# row 7
$ws.Cells.Item(7, 5).Value = 0.51
$ws.Cells.Item(7, 14).Value = 0.643
$ws.Cells.Item(7, 17).Value = 0.016
$ws.Cells.Item(7, 23).Value = 0.275
$ws.Cells.Item(7, 35).Value = 0.305
$ws.Cells.Item(7, 47).Value = 0.307
$ws.Cells.Item(7, 53).Value = 1.429
$ws.Cells.Item(7, 59).Value = 0.438
$ws.Cells.Item(7, 65).Value = 0.576
$ws.Cells.Item(7, 68).Value = 0.476
$ws.Cells.Item(7, 69).Value = 0.49

# row 8
$ws.Cells.Item(8, 5).Value = 0.605
$ws.Cells.Item(8, 6).Value = 0.112
$ws.Cells.Item(8, 7).Value = 0.334
$ws.Cells.Item(8, 14).Value = 0.773
$ws.Cells.Item(8, 15).Value = 0.066
$ws.Cells.Item(8, 16).Value = 0.256
$ws.Cells.Item(8, 17).Value = 0.017
$ws.Cells.Item(8, 19).Value = 0.107
$ws.Cells.Item(8, 23).Value = 0.301
$ws.Cells.Item(8, 25).Value = 0.347
$ws.Cells.Item(8, 35).Value = 0.33
$ws.Cells.Item(8, 36).Value = 0.129
$ws.Cells.Item(8, 37).Value = 0.359
$ws.Cells.Item(8, 47).Value = 0.308
$ws.Cells.Item(8, 48).Value = 0.08500000000000001
$ws.Cells.Item(8, 49).Value = 0.291
$ws.Cells.Item(8, 53).Value = 1.75
$ws.Cells.Item(8, 54).Value = 0.124
$ws.Cells.Item(8, 55).Value = 0.352
$ws.Cells.Item(8, 59).Value = 0.5679999999999999
$ws.Cells.Item(8, 60).Value = 0.106
$ws.Cells.Item(8, 61).Value = 0.326
$ws.Cells.Item(8, 65).Value = 0.6909999999999999
$ws.Cells.Item(8, 66).Value = 0.066
$ws.Cells.Item(8, 67).Value = 0.258
$ws.Cells.Item(8, 68).Value = 0.583
$ws.Cells.Item(8, 69).Value = 0.605

# row 9
$ws.Cells.Item(9, 5).Value = 0.553
$ws.Cells.Item(9, 6).Value = 0.247
$ws.Cells.Item(9, 7).Value = 0.497
$ws.Cells.Item(9, 14).Value = 0.681
$ws.Cells.Item(9, 15).Value = 0.217
$ws.Cells.Item(9, 16).Value = 0.466
$ws.Cells.Item(9, 23).Value = 0.202
$ws.Cells.Item(9, 24).Value = 0.161
$ws.Cells.Item(9, 25).Value = 0.402
$ws.Cells.Item(9, 35).Value = 0.255
$ws.Cells.Item(9, 36).Value = 0.19
$ws.Cells.Item(9, 37).Value = 0.436
$ws.Cells.Item(9, 53).Value = 1.712
$ws.Cells.Item(9, 54).Value = 0.248
$ws.Cells.Item(9, 55).Value = 0.498
$ws.Cells.Item(9, 59).Value = 0.606
$ws.Cells.Item(9, 60).Value = 0.239
$ws.Cells.Item(9, 61).Value = 0.489
$ws.Cells.Item(9, 65).Value = 0.649
$ws.Cells.Item(9, 66).Value = 0.228
$ws.Cells.Item(9, 67).Value = 0.477
$ws.Cells.Item(9, 68).Value = 0.571
$ws.Cells.Item(9, 69).Value = 0.588

# row 10
$ws.Cells.Item(10, 5).Value = 0.681
$ws.Cells.Item(10, 6).Value = 0.217
$ws.Cells.Item(10, 7).Value = 0.466
$ws.Cells.Item(10, 14).Value = 0.872
$ws.Cells.Item(10, 15).Value = 0.111
$ws.Cells.Item(10, 16).Value = 0.334
$ws.Cells.Item(10, 23).Value = 0.372
$ws.Cells.Item(10, 24).Value = 0.234
$ws.Cells.Item(10, 25).Value = 0.483
$ws.Cells.Item(10, 35).Value = 0.362
$ws.Cells.Item(10, 36).Value = 0.231
$ws.Cells.Item(10, 37).Value = 0.48
$ws.Cells.Item(10, 47).Value = 0.298
$ws.Cells.Item(10, 48).Value = 0.209
$ws.Cells.Item(10, 49).Value = 0.457
$ws.Cells.Item(10, 53).Value = 2.085
$ws.Cells.Item(10, 54).Value = 0.243
$ws.Cells.Item(10, 55).Value = 0.493
$ws.Cells.Item(10, 59).Value = 0.66
$ws.Cells.Item(10, 60).Value = 0.225
$ws.Cells.Item(10, 61).Value = 0.474
$ws.Cells.Item(10, 65).Value = 0.84
$ws.Cells.Item(10, 66).Value = 0.134
$ws.Cells.Item(10, 67).Value = 0.366
$ws.Cells.Item(10, 68).Value = 0.695
$ws.Cells.Item(10, 69).Value = 0.726

# row 11
$ws.Cells.Item(11, 5).Value = 0.713
$ws.Cells.Item(11, 6).Value = 0.205
$ws.Cells.Item(11, 7).Value = 0.452
$ws.Cells.Item(11, 14).Value = 0.894
$ws.Cells.Item(11, 15).Value = 0.095
$ws.Cells.Item(11, 16).Value = 0.308
$ws.Cells.Item(11, 23).Value = 0.372
$ws.Cells.Item(11, 24).Value = 0.234
$ws.Cells.Item(11, 25).Value = 0.483
$ws.Cells.Item(11, 35).Value = 0.394
$ws.Cells.Item(11, 36).Value = 0.239
$ws.Cells.Item(11, 37).Value = 0.489
$ws.Cells.Item(11, 47).Value = 0.436
$ws.Cells.Item(11, 48).Value = 0.246
$ws.Cells.Item(11, 49).Value = 0.496
$ws.Cells.Item(11, 53).Value = 2.085
$ws.Cells.Item(11, 54).Value = 0.243
$ws.Cells.Item(11, 55).Value = 0.493
$ws.Cells.Item(11, 59).Value = 0.66
$ws.Cells.Item(11, 60).Value = 0.225
$ws.Cells.Item(11, 61).Value = 0.474
$ws.Cells.Item(11, 65).Value = 0.84
$ws.Cells.Item(11, 66).Value = 0.134
$ws.Cells.Item(11, 67).Value = 0.366
$ws.Cells.Item(11, 68).Value = 0.695
$ws.Cells.Item(11, 69).Value = 0.728

# row 12
$ws.Cells.Item(12, 5).Value = 1.403
$ws.Cells.Item(12, 6).Value = 0.748
$ws.Cells.Item(12, 7).Value = 0.865
$ws.Cells.Item(12, 14).Value = 1.465
$ws.Cells.Item(12, 15).Value = 1.039
$ws.Cells.Item(12, 16).Value = 1.02
$ws.Cells.Item(12, 23).Value = 1.629
$ws.Cells.Item(12, 24).Value = 0.576
$ws.Cells.Item(12, 25).Value = 0.759
$ws.Cells.Item(12, 35).Value = 1.703
$ws.Cells.Item(12, 36).Value = 1.29
$ws.Cells.Item(12, 37).Value = 1.136
$ws.Cells.Item(12, 47).Value = 2.767
$ws.Cells.Item(12, 48).Value = 2.737
$ws.Cells.Item(12, 49).Value = 1.654
$ws.Cells.Item(12, 54).Value = 0.395
$ws.Cells.Item(12, 55).Value = 0.629
$ws.Cells.Item(12, 59).Value = 1.097
$ws.Cells.Item(12, 60).Value = 0.12
$ws.Cells.Item(12, 61).Value = 0.346
$ws.Cells.Item(12, 65).Value = 1.291
$ws.Cells.Item(12, 66).Value = 0.333
$ws.Cells.Item(12, 67).Value = 0.577
$ws.Cells.Item(12, 68).Value = 1.232
$ws.Cells.Item(12, 69).Value = 1.257

# row 13
$ws.Cells.Item(13, 5).Value = 1.573
$ws.Cells.Item(13, 6).Value = 0.652
$ws.Cells.Item(13, 7).Value = 0.8080000000000001
$ws.Cells.Item(13, 14).Value = 2.068
$ws.Cells.Item(13, 15).Value = 0.928
$ws.Cells.Item(13, 16).Value = 0.964
$ws.Cells.Item(13, 23).Value = 1.037
$ws.Cells.Item(13, 24).Value = 0.193
$ws.Cells.Item(13, 25).Value = 0.439
$ws.Cells.Item(13, 35).Value = 1.28
$ws.Cells.Item(13, 36).Value = 0.37
$ws.Cells.Item(13, 37).Value = 0.608
$ws.Cells.Item(13, 47).Value = 2.285
$ws.Cells.Item(13, 48).Value = 0.925
$ws.Cells.Item(13, 49).Value = 0.962
$ws.Cells.Item(13, 53).Value = 2.344
$ws.Cells.Item(13, 54).Value = 0.296
$ws.Cells.Item(13, 55).Value = 0.544
$ws.Cells.Item(13, 59).Value = 0.584
$ws.Cells.Item(13, 60).Value = 0.07099999999999999
$ws.Cells.Item(13, 61).Value = 0.266
$ws.Cells.Item(13, 65).Value = 0.892
$ws.Cells.Item(13, 66).Value = 0.282
$ws.Cells.Item(13, 67).Value = 0.531
$ws.Cells.Item(13, 68).Value = 0.781
$ws.Cells.Item(13, 69).Value = 0.725
